$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values: PRB/Perabot -> BJU/Baju
$ws.Range("A2").Value = "BJU"
$ws.Range("B2").Value = "Baju"

# Remove row 3 entirely (ALTMD / Alat Mandi)
$ws.Rows(3).Delete()

# Update selection to B2 as in the saved file
$ws.Range("B2").Select()
